# Generate Report for Handback
# Updates the localization-status workbook to reflect a handback transform
# failure for the d85be330-... file in both the zh-cn and de-de locales,
# and records the corresponding error detail message.

$wb = $excel.ActiveWorkbook

$statusFailed = "Handback transform failed"

$zhErrorDetail = "Handback file name: rxnvzz0t.xf2 is different with handoff file name: d85be330-a85a-4c0e-9a28-e105fa3d66f1.a788c72ebeff1b1e564ccadd5652ac7a7c8b25cd.zh-cn."
$deErrorDetail = "Handback file name: rxnvzz0t.xf2 is different with handoff file name: d85be330-a85a-4c0e-9a28-e105fa3d66f1.a788c72ebeff1b1e564ccadd5652ac7a7c8b25cd.de-de."

# Overview sheet: row for d85be330-...md (row 3) shows the status per locale
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $statusFailed
$wsOverview.Range("C3").Value = $statusFailed

# zh-cn sheet: row 3 corresponds to the d85be330-... handoff entry
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $statusFailed
$wsZhCn.Range("K3").Value = $zhErrorDetail

# de-de sheet: row 3 corresponds to the d85be330-... handoff entry
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $statusFailed
$wsDeDe.Range("K3").Value = $deErrorDetail
